$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Mine" (J) column count for the Uncommon rarity row (row 9)
$ws.Range("J9").Value = 26

# Update the "Mine" (E) column count for the Cost=2 row (row 10)
$ws.Range("E10").Value = 15

# Update the "Mine" (J) column count for the Skill type row (row 18)
$ws.Range("J18").Value = 29

# Move the active selection to J19, as in the author's edit
$ws.Range("J19").Select()
